$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 2: "Program Easiest Level" (same text already in E2) gets copied across F2:H2
$progEasiest = $ws.Range("E2").Value()
$ws.Range("F2").Value = $progEasiest
$ws.Range("G2").Value = $progEasiest
$ws.Range("H2").Value = $progEasiest

# New progress notes typed into column H (authoring order matches shared-string order)
$ws.Range("H5").Value = "Plan for this week is levels"
$ws.Range("H6").Value = "Two main mechanics id'd: Jumping and block pushing"
$ws.Range("H3").Value = "Write some level ""stories"" similar to what Maddy Thorsen was talking about"
$ws.Range("H4").Value = "Brainstorm Ideas"
$ws.Range("H7").Value = "Maybe attempt Thorsen's story exercise on levels from platformers you like?"
$ws.Range("H8").Value = "Watch Mark Brown's level design videos?"

# H6 gets the same yellow-highlight / centered / wrap-text look used elsewhere in the workbook
$ws.Range("H6").Interior.ColorIndex = 6
$ws.Range("H6").HorizontalAlignment = -4108
$ws.Range("H6").WrapText = $true

# Column H has to widen considerably to fit the new, much longer text
$ws.Columns.Item(8).ColumnWidth = 68.75

# Selection / scroll position moved further down and right as notes were added
$ws.Range("H9").Select()
